# The presentation's custom "Integral" theme colour scheme (applied to the
# single slide master / theme1.xml) is swapped out for the built-in
# "Office" colour palette (the palette that currently lives in theme2.xml,
# the presentation's secondary/notes theme part).
#
# PowerPoint's ColorScheme object exposes the 12 DrawingML theme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) as Colors(1..12); writing
# to them updates the <a:clrScheme> of the theme part backing the active
# slide master (theme1.xml).

function RgbToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Office theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = RgbToOle $officeColors[$i - 1]
}
